# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped values (per the GitHub Actions refresh commit).
# All of these cells hold plain text in the workbook (not numbers), so for
# any new value that looks numeric we prefix it with a leading apostrophe
# (the same trick Excel itself uses) to force it to stay text instead of
# being auto-converted to a Number cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.756.63'
$ws.Range('E2').Value = '  +1.54%  '
$ws.Range('D3').Value = '2.802.63'
$ws.Range('E3').Value = '  +1.83%  '
$ws.Range('D5').Value = '''353.04'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').Value = '''112.12'
$ws.Range('E6').Value = '  +4.70%  '
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.622'
$ws.Range('E9').Value = '  +7.61%  '
$ws.Range('D10').Value = '''40.15'
$ws.Range('E10').Value = '  +2.97%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').Value = '''0.0838'
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('D13').Value = '''19.91'
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').Value = '''7.76'
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('D15').Value = '3.243.51'
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('D16').Value = '2.818.93'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('D17').Value = '''0.956'
$ws.Range('E17').Value = '  +3.64%  '
$ws.Range('D18').Value = '51.780.74'
$ws.Range('E18').Value = '  +1.57%  '
$ws.Range('D19').Value = '''7.62'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').Value = '''3.29'
$ws.Range('E20').Value = '  +8.76%  '
$ws.Range('E21').Value = '  +4.80%  '
$ws.Range('D22').Value = '0.0₃0973'
$ws.Range('E22').Value = '  +2.12%  '
$ws.Range('D23').Value = '''70.25'
$ws.Range('E23').Value = '  +1.78%  '
$ws.Range('D24').Value = '''267.56'
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D27').Value = '''26.13'
$ws.Range('E27').Value = '  +1.39%  '
$ws.Range('E28').Value = '  +0.95%  '
$ws.Range('D29').Value = '''39.14'
$ws.Range('E29').Value = '  +14.10%  '
$ws.Range('E30').Value = '  +3.93%  '
$ws.Range('E31').Value = '  +2.25%  '
$ws.Range('D32').Value = '''52.27'
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  +1.99%  '
$ws.Range('D34').Value = '''0.0900'
$ws.Range('E34').Value = '  +9.07%  '
$ws.Range('E35').Value = '  +2.63%  '
$ws.Range('D36').Value = '''5.52'
$ws.Range('E36').Value = '  +4.85%  '
$ws.Range('E37').Value = '  -0.05%  '
$ws.Range('E38').Value = '  +4.45%  '
$ws.Range('E39').Value = '  +1.31%  '
$ws.Range('E41').Value = '  +2.41%  '
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('D44').Value = '''119.78'
$ws.Range('E44').Value = '  -0.70%  '
$ws.Range('D45').Value = '''21.94'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '''3.52'
$ws.Range('E46').Value = '  +10.16%  '
$ws.Range('D47').Value = '''2.49'
$ws.Range('E47').Value = '  +9.58%  '
$ws.Range('D48').Value = '2.118.56'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').Value = '''0.977'
$ws.Range('E49').Value = '  +8.08%  '
$ws.Range('D50').Value = '''5.48'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('E51').Value = '  +8.24%  '
